$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 177: LeetCode 1488 - Avoid Flood in The City ---
$ws.Range("A177").Value = 1488
$ws.Range("B177").Value = "Avoid Flood in The City"
$ws.Range("C177").Value = "#greedy #binary-search "
$ws.Range("D177").Value = "medium"
$ws.Range("E177").Value = 0
$ws.Range("F177").Value = 1
$ws.Range("G177").Value = 20
$ws.Range("H177").Value = 45937
$ws.Range("I177").Value = 45937

# --- Row 178: LeetCode 2300 - Successful Pairs of Spells and Potions ---
$ws.Range("A178").Value = 2300
$ws.Range("B178").Value = "Successful Pairs of Spells and Potions"
$ws.Range("C178").Value = "#binary-search #sorting "
$ws.Range("D178").Value = "medium"
$ws.Range("E178").Value = 1
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 21
$ws.Range("H178").Value = 45938
$ws.Range("I178").Value = 45938

# --- Row 179: LeetCode 774 - Minimize Max Distance to Gas Station ---
$ws.Range("A179").Value = 774
$ws.Range("B179").Value = "Minimize Max Distance to Gas Station"
$ws.Range("C179").Value = "#binary-search"
$ws.Range("D179").Value = "hard"
$ws.Range("E179").Value = 1
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 20
$ws.Range("H179").Value = 45938
$ws.Range("I179").Value = 45938

# Dates in columns H/I need the same date number format / cell style as the
# rest of the table (col default style is generic, the date look comes from
# a per-cell style applied on the existing rows) - copy it down from row 176.
$ws.Range("H176:I176").Copy()
$ws.Range("H177:I179").PasteSpecial(-4122)

# Match the row heights the sheet ends up with once the new rows hold
# (wrapped) text, same as every other data row in the table.
$ws.Rows.Item(177).RowHeight = 34
$ws.Rows.Item(178).RowHeight = 34
$ws.Rows.Item(179).RowHeight = 34

# Leave the selection where Excel would after the last edit (the final cell
# touched, H179:I179) - matches the sheetView state after data entry.
$null = $ws.Range("H179:I179").Select()
